$d = $word.ActiveDocument

# The document contains four paragraphs of the form:
#   <run1: "<id>" (Courier New, color 7f6000)>
#   <run2: "p158v_N" (color 000000)>
#   <run3: "</id>" (Courier New, color 7f6000)>
# which need to be merged into a single run:
#   <run: "<id>p158v_N</id>">
# using the formatting of the first run (Courier New / 7f6000), for
# N = 1, 2, 3, 4.
#
# Find & Replace across the run boundary collapses the matched text into
# a single run that inherits the formatting of the first run in the
# match, which is exactly the merge we need.

for ($i = 1; $i -le 4; $i++) {
    $needle = "<id>p158v_$i</id>"
    $d.Content.Find.Execute($needle, $false, $false, $false, $false, $false, $true, 1, $false, $needle, 2)
}
